$wb = $excel.ActiveWorkbook

# Rename sheet "Model3" to "Model310"
$wsModel3 = $wb.Worksheets.Item("Model3")
$wsModel3.Name = "Model310"

# --- Update Model1 sheet (column D: solution time) ---
$wsModel1 = $wb.Worksheets.Item("Model1")
$wsModel1.Range("D2").Value = 0.0150001049041748
$wsModel1.Range("D3").Value = 0.01900005340576172
$wsModel1.Range("D4").Value = 0.01399993896484375
$wsModel1.Range("D5").Value = 0.01399993896484375
$wsModel1.Range("D6").Value = 0.01799988746643066
$wsModel1.Range("D7").Value = 0.01699995994567871
$wsModel1.Range("D8").Value = 0.01100015640258789
$wsModel1.Range("D10").Value = 0.01999998092651367
$wsModel1.Range("D11").Value = 0.0130000114440918
$wsModel1.Range("D12").Value = 0.01600003242492676
$wsModel1.Range("D13").Value = 0.0130000114440918
$wsModel1.Range("D14").Value = 0.0130000114440918
$wsModel1.Range("D15").Value = 0.0149998664855957
$wsModel1.Range("D16").Value = 0.0130000114440918
$wsModel1.Range("D17").Value = 0.01699995994567871
$wsModel1.Range("D18").Value = 0.01200008392333984
$wsModel1.Range("D19").Value = 0.01699995994567871
$wsModel1.Range("D20").Value = 0.01200008392333984
$wsModel1.Range("D21").Value = 0.0130000114440918

# --- Update Model2 sheet (column D: solution time) ---
$wsModel2 = $wb.Worksheets.Item("Model2")
$wsModel2.Range("D2").Value = 0.4739999771118164
$wsModel2.Range("D3").Value = 0.9619998931884766
$wsModel2.Range("D4").Value = 0.2660000324249268
$wsModel2.Range("D5").Value = 0.1120002269744873
$wsModel2.Range("D6").Value = 0.1119999885559082
$wsModel2.Range("D7").Value = 0.1150000095367432
$wsModel2.Range("D8").Value = 0.07299995422363281
$wsModel2.Range("D9").Value = 0.03099989891052246
$wsModel2.Range("D10").Value = 0.1460001468658447
$wsModel2.Range("D11").Value = 0.1289999485015869
$wsModel2.Range("D12").Value = 0.1559998989105225
$wsModel2.Range("D13").Value = 0.1390001773834229
$wsModel2.Range("D14").Value = 0.18399977684021
$wsModel2.Range("D15").Value = 0.2349998950958252
$wsModel2.Range("D16").Value = 0.1140000820159912
$wsModel2.Range("D17").Value = 0.2019999027252197
$wsModel2.Range("D18").Value = 0.1150000095367432
$wsModel2.Range("D19").Value = 0.08299994468688965
$wsModel2.Range("D20").Value = 0.1400001049041748
$wsModel2.Range("D21").Value = 0.0280001163482666

# --- Update Model310 sheet (columns C, D, E, F, G) ---
$wsModel310 = $wb.Worksheets.Item("Model310")
$wsModel310.Range("C2").Value = 10770.99999469879
$wsModel310.Range("D2").Value = 3600.618000030518
$wsModel310.Range("E2").Value = 276
$wsModel310.Range("G2").Value = 227752
$wsModel310.Range("C3").Value = 12494.99997567427
$wsModel310.Range("D3").Value = 3601.017999887466
$wsModel310.Range("E3").Value = 308
$wsModel310.Range("G3").Value = 279398
$wsModel310.Range("C4").Value = 9518.999783110961
$wsModel310.Range("D4").Value = 3600.436000108719
$wsModel310.Range("E4").Value = 258
$wsModel310.Range("G4").Value = 201072
$wsModel310.Range("C5").Value = 6598.999934106854
$wsModel310.Range("D5").Value = 3600.695000171661
$wsModel310.Range("E5").Value = 190
$wsModel310.Range("G5").Value = 99268
$wsModel310.Range("C6").Value = 4997.999683341697
$wsModel310.Range("D6").Value = 3600.507999897003
$wsModel310.Range("E6").Value = 166
$wsModel310.Range("F6").Value = "maxTimeLimit"
$wsModel310.Range("G6").Value = 69916
$wsModel310.Range("C7").Value = 5720.999254496943
$wsModel310.Range("D7").Value = 3600.757999897003
$wsModel310.Range("E7").Value = 174
$wsModel310.Range("F7").Value = "maxTimeLimit"
$wsModel310.Range("G7").Value = 79964
$wsModel310.Range("C8").Value = 4134.999177478211
$wsModel310.Range("D8").Value = 3607.441999912262
$wsModel310.Range("E8").Value = 156
$wsModel310.Range("F8").Value = "maxTimeLimit"
$wsModel310.Range("G8").Value = 66950
$wsModel310.Range("C9").Value = 1601.999998358191
$wsModel310.Range("D9").Value = 0.8029999732971191
$wsModel310.Range("E9").Value = 25
$wsModel310.Range("G9").Value = 3706
$wsModel310.Range("C10").Value = 8575.999372841841
$wsModel310.Range("D10").Value = 3600.623000144958
$wsModel310.Range("E10").Value = 250
$wsModel310.Range("G10").Value = 176244
$wsModel310.Range("C11").Value = 6138.999382766653
$wsModel310.Range("D11").Value = 3600.554000139236
$wsModel310.Range("E11").Value = 210
$wsModel310.Range("F11").Value = "maxTimeLimit"
$wsModel310.Range("G11").Value = 125484
$wsModel310.Range("C12").Value = 7148.999629195149
$wsModel310.Range("D12").Value = 3600.576999902725
$wsModel310.Range("E12").Value = 220
$wsModel310.Range("G12").Value = 135766
$wsModel310.Range("C13").Value = 5313.513460910006
$wsModel310.Range("D13").Value = 3634.87700009346
$wsModel310.Range("E13").Value = 192
$wsModel310.Range("F13").Value = "maxTimeLimit"
$wsModel310.Range("G13").Value = 102906
$wsModel310.Range("C14").Value = 7108.998884280798
$wsModel310.Range("D14").Value = 3600.973999977112
$wsModel310.Range("E14").Value = 220
$wsModel310.Range("F14").Value = "maxTimeLimit"
$wsModel310.Range("G14").Value = 141346
$wsModel310.Range("C15").Value = 8092.999516291501
$wsModel310.Range("D15").Value = 3600.71799993515
$wsModel310.Range("E15").Value = 240
$wsModel310.Range("G15").Value = 166110
$wsModel310.Range("C16").Value = 7540.621403651022
$wsModel310.Range("D16").Value = 3600.673000097275
$wsModel310.Range("E16").Value = 218
$wsModel310.Range("G16").Value = 135604
$wsModel310.Range("C17").Value = 5816.999102917282
$wsModel310.Range("D17").Value = 3600.576999902725
$wsModel310.Range("E17").Value = 192
$wsModel310.Range("F17").Value = "maxTimeLimit"
$wsModel310.Range("G17").Value = 103886
$wsModel310.Range("C18").Value = 4574.998437223135
$wsModel310.Range("D18").Value = 3600.417000055313
$wsModel310.Range("E18").Value = 156
$wsModel310.Range("F18").Value = "maxTimeLimit"
$wsModel310.Range("G18").Value = 66950
$wsModel310.Range("C19").Value = 3615
$wsModel310.Range("D19").Value = 2692.256999969482
$wsModel310.Range("E19").Value = 134
$wsModel310.Range("G19").Value = 44894
$wsModel310.Range("C20").Value = 6505.999870973515
$wsModel310.Range("D20").Value = 3600.544999837875
$wsModel310.Range("E20").Value = 224
$wsModel310.Range("F20").Value = "maxTimeLimit"
$wsModel310.Range("G20").Value = 142182
$wsModel310.Range("D21").Value = 0.1380000114440918
$wsModel310.Range("E21").Value = 7
$wsModel310.Range("G21").Value = 154
